# Automatische test-sync: 2025-08-03 14:51:50
# Append the new "Testmail #9" log row to the Logs sheet and refresh the
# Dashboard category summary table to match.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append row 19 -------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A19").Value = "Hoi, hebben jullie al iets gehoord?"
$logs.Range("B19").Value = "mailmind.test@zohomail.eu"
$logs.Range("C19").Value = "Testmail #9: Hoi, hebben jullie al iets gehoord?"
$logs.Range("D19").Value = "Overig"
$logs.Range("E19").Value = "Beste klant,`nBedankt voor uw e-mail. Kunt u ons meer informatie geven over waar u precies op wacht of waar u over gehoord wilt hebben? Zo kunnen wij u beter van dienst zijn.`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$logs.Range("F19").Value = "2025-08-03 14:51:12"
$logs.Range("G19").Value = "Ja"
$logs.Range("H19").Value = "Nee"
$logs.Range("I19").Value = "Ja"
$logs.Range("J19").Value = "Nee"

# --- Logs sheet: extend conditional formatting ranges to include row 19 --
foreach ($col in @("D", "G", "H", "I", "J")) {
    $newRange = $logs.Range("$col" + "2:" + "$col" + "19")
    $fc = $newRange.FormatConditions.Item(1)
    $fc.ModifyAppliesToRange($newRange)
}

# --- Dashboard sheet: refresh the category summary ------------------------
# New row 19 is "Overig", so its count goes from 4 -> 5, and "Planning /
# Afspraak" drops to the row below it in the summary ordering.
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A3").Value = "Overig"
$dash.Range("B3").Value = 5
$dash.Range("A4").Value = "Planning / Afspraak"
$dash.Range("B4").Value = 4
